# Update a handful of numeric values on the (single) worksheet to reflect
# a refreshed RandomForest imputation run ("Update Name of Algo").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = 6.006799999999996
$ws.Range("A3").Value  = -21.48760000000003
$ws.Range("B5").Value  = 4.730600000000003
$ws.Range("C5").Value  = -14.0178
$ws.Range("D7").Value  = -7.189099999999995
$ws.Range("C9").Value  = -12.08410000000001
$ws.Range("C11").Value = -12.7745
$ws.Range("D11").Value = -8.116700000000002
$ws.Range("A14").Value = -20.39389999999998
$ws.Range("D19").Value = -8.511899999999997
$ws.Range("A21").Value = -21.08030000000001
$ws.Range("C21").Value = -10.71909999999999
$ws.Range("D21").Value = -7.2479
$ws.Range("A23").Value = -21.46660000000002
$ws.Range("A25").Value = -22.65460000000003
